$wb = $excel.ActiveWorkbook

$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# VENTAS POR GRUPO: M19 (PORCELANATO for PUEBLA GONZALEZ MARIO DANIEL) 0 -> 3252.41
$wsVentasGrupo.Range("M19").Value = 3252.41
# VENTAS POR GRUPO: M25 "4 de 23" -> "5 de 23"
$wsVentasGrupo.Range("M25").Value = "5 de 23"

# VENTA MENSUAL: F19 (septiembre for PUEBLA GONZALEZ MARIO DANIEL) 0 -> 3252.41
$wsVentaMensual.Range("F19").Value = 3252.41
# VENTA MENSUAL: F25 total septiembre 30206.69 -> 33459.1
$wsVentaMensual.Range("F25").Value = 33459.1

# CUMPLIMIENTO MENSUAL: PORCELANATO row 12
$wsCumplimiento.Range("D12").Value = 32645.98
$wsCumplimiento.Range("E12").Value = 10454.1054117774
$wsCumplimiento.Range("F12").Value = 0.757445830747223

# CUMPLIMIENTO MENSUAL: TOTAL row 15
$wsCumplimiento.Range("D15").Value = 33459.1
$wsCumplimiento.Range("E15").Value = 24744.36623249458
$wsCumplimiento.Range("F15").Value = 0.5748643880820971
